$d = $word.ActiveDocument
try { $d.UpdateStylesOnOpen = $true; Write-Output "set ok" } catch { Write-Output "ERR: $_" }
try { Write-Output $d.AttachedTemplate } catch { Write-Output "ERR2: $_" }
